# Filter - Study - Test Suit
# The "startup" sheet's Cases row actually drives the Participants query,
# so rename the TabName value from "CasesTab" to "ParticipantsTab".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ParticipantsTab"

# Move/restore the active selection to A2 (matches the saved view state).
$ws.Range("A2").Select()
